# "datos del dia 4-04" - fill in the KC daily report table with the
# records collected for 4-04 (rows 2-8), extend the table / validations
# to match, and update dimension + selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Enter the raw values first (while rows 3-8 are still blank /
#        General-formatted) so numeric literals land as real numbers and
#        not as text - formatting is copied down afterwards. --------------

# Row 2
$ws.Range("A2").Value = 45021
$ws.Range("B2").Value = "9651"
$ws.Range("C2").Value = "Madolche"
$ws.Range("D2").Value = "Robar sentido Nivel bajo"
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0
$ws.Range("H2").Value = 1
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 1

# Row 3
$ws.Range("A3").Value = 45021
$ws.Range("B3").Value = 9881
$ws.Range("C3").Value = "Rokket"
$ws.Range("D3").Value = "Lanzamiento Borre"
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 0

# Row 4
$ws.Range("A4").Value = 45021
$ws.Range("B4").Value = 2947
$ws.Range("C4").Value = "Rose Dragon"
$ws.Range("D4").Value = "Conjuro de Rosas"
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 0

# Row 5
$ws.Range("A5").Value = 45021
$ws.Range("B5").Value = 2848
$ws.Range("C5").Value = "Solfachord"
$ws.Range("D5").Value = "Péndulos Unidos"
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 0

# Row 6
$ws.Range("A6").Value = 45021
$ws.Range("B6").Value = 7906
$ws.Range("C6").Value = "Synchrons"
$ws.Range("D6").Value = "Un Vínculo Ilumina el Futuro"
$ws.Range("F6").Value = 0
$ws.Range("G6").Value = 1
$ws.Range("H6").Value = 0
$ws.Range("I6").Value = 0
$ws.Range("J6").Value = 0

# Row 7
$ws.Range("A7").Value = 45021
$ws.Range("B7").Value = 770
$ws.Range("C7").Value = "Infernoid"
$ws.Range("D7").Value = "Aumento de LP alfa"
$ws.Range("F7").Value = 0
$ws.Range("G7").Value = 1
$ws.Range("H7").Value = 1
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 0

# Row 8 (note: D8 is entered before C8 so the shared-string table picks
# up "Alma Rugiente" before "Resonators", matching the author's original
# entry order)
$ws.Range("A8").Value = 45021
$ws.Range("B8").Value = 3338
$ws.Range("D8").Value = "Alma Rugiente"
$ws.Range("C8").Value = "Resonators"
$ws.Range("F8").Value = 0
$ws.Range("G8").Value = 0
$ws.Range("H8").Value = 1
$ws.Range("I8").Value = 0
$ws.Range("J8").Value = 0

# --- 2. Copy the formatting already present on row 2 down to rows 3-8 so
#        every new row picks up the same number formats / fonts (styles
#        2,6,1,1,3,4,1,1,1,1,1 for columns A..K) without minting new xf
#        records or turning the numbers already entered into text. --------
$ws.Range("A2:K2").Copy() | Out-Null
$ws.Range("A3:K8").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# --- 3. Calculated columns E (ndmax) and K (day) for the new rows --------
$ws.Range("E3").Formula = "=IF(`$A3=`"`",`"`",`$A3-1)"
$ws.Range("E4").Formula = "=IF(`$A4=`"`",`"`",`$A4-1)"
$ws.Range("E5").Formula = "=IF(`$A5=`"`",`"`",`$A5-1)"
$ws.Range("E6").Formula = "=IF(`$A6=`"`",`"`",`$A6-1)"
$ws.Range("E7").Formula = "=IF(`$A7=`"`",`"`",`$A7-1)"
$ws.Range("E8").Formula = "=IF(`$A8=`"`",`"`",`$A8-1)"

$ws.Range("K3").Formula = "=IFERROR(ROUND((VALUE(TEXT(`$E3,`"DD`"))),0),`"`")"
$ws.Range("K4").Formula = "=IFERROR(ROUND((VALUE(TEXT(`$E4,`"DD`"))),0),`"`")"
$ws.Range("K5").Formula = "=IFERROR(ROUND((VALUE(TEXT(`$E5,`"DD`"))),0),`"`")"
$ws.Range("K6").Formula = "=IFERROR(ROUND((VALUE(TEXT(`$E6,`"DD`"))),0),`"`")"
$ws.Range("K7").Formula = "=IFERROR(ROUND((VALUE(TEXT(`$E7,`"DD`"))),0),`"`")"
$ws.Range("K8").Formula = "=IFERROR(ROUND((VALUE(TEXT(`$E8,`"DD`"))),0),`"`")"

# --- 4. Expand the table / autofilter to cover the new rows --------------
$tbl = $ws.ListObjects.Item(1)
$tbl.Resize($ws.Range("A1:K8"))

# --- 5. Expand the data validations that were scoped to row 2 only -------
$ws.Range("F2:J2").Validation.Delete()
$ws.Range("F2:J8").Validation.Add(1, 1, 1, "0", "1")
$ws.Range("F2:J8").Validation.ErrorTitle = "Error de Tipeo"
$ws.Range("F2:J8").Validation.ErrorMessage = "solo es 1 u 0 para definir v o f"

$ws.Range("B2").Validation.Delete()
$ws.Range("B2:B8").Validation.Add(7, 1, 1, "COUNTIF(`$B`$2:`$B`$668,`$B2)=1")
$ws.Range("B2:B8").Validation.ErrorTitle = "Usuario Existente"
$ws.Range("B2:B8").Validation.ErrorMessage = "corrija el usuario o verifique si son correctos los datos"
$ws.Range("B2:B8").Validation.InCellDropdown = $false

# --- 6. Selection moves to G6 after data entry ----------------------------
$ws.Range("G6").Select()

Write-Output "done"
